$wb = $excel.ActiveWorkbook

# --- Sheet1: remove the stored login credentials (rows 1-2: "User Name"/"Password"
#     header row and the actual "AgsautoT04"/"SERVICE$08" values) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1:B2").ClearContents()

# --- Sheet2: bump the schedule-number counter ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B1").Value = 4

# --- Make Sheet1 the active/visible sheet (it was Sheet2 before) and reset
#     the selection back to the top-left of the cleared area ---
$ws1.Activate()
$ws1.Range("A1:B2").Select()
